# Auto-applied price/volume refresh for cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.516.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "'2.317.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'518.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").Value = "'135.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.79%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").Value = "'2.332.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("D10").Value = "'0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "'5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.62%  "
$ws.Range("D13").Value = "'0.342"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "'24.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "'2.724.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "'56.574.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("D18").Value = "'2.335.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("D19").Value = "'10.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "'4.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'322.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.06%  "
$ws.Range("D22").Value = "'6.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'60.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("E25").Value = "  +5.95%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'7.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.97%  "
$ws.Range("D28").Value = "'1.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.93%  "
$ws.Range("D29").Value = "'0.0₃0741"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.49%  "
$ws.Range("E30").Value = "  +4.60%  "
$ws.Range("D31").Value = "'166.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'6.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").Value = "'18.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'0.991"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "'1.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("D37").Value = "'0.925"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("D38").Value = "'4.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("D39").Value = "'1.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.51%  "
$ws.Range("D40").Value = "'37.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").Value = "'0.381"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").Value = "'140.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.40%  "
$ws.Range("D43").Value = "'3.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.41%  "
$ws.Range("D44").Value = "'5.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.93%  "
$ws.Range("D45").Value = "'278.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.26%  "
$ws.Range("D46").Value = "'0.0933"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("D47").Value = "'0.0508"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "'0.561"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("D49").Value = "'0.0218"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("D50").Value = "'0.380"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "'17.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.70%  "
